# Apply weekly price-update edit: insert one new data row at row 1171
# (pushing the existing rows 1171:1255 down to 1172:1256), and populate the
# new row with the latest "Packham's Triumph" / "Primera" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 1171. Excel shifts rows 1171:1255
# down to 1172:1256 and copies formatting from the row above (row 1170),
# which already carries the date number format used in column D.
$ws.Rows.Item(1171).Insert()

# Populate the newly inserted row 1171 with the new observation.
$ws.Cells.Item(1171, 1).Value = 5
$ws.Cells.Item(1171, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1171, 3).Value = "Maule"
$ws.Cells.Item(1171, 4).Value = 45265
$ws.Cells.Item(1171, 5).Value = 7
$ws.Cells.Item(1171, 6).Value = "Fruta"
$ws.Cells.Item(1171, 7).Value = 100104
$ws.Cells.Item(1171, 8).Value = "Frutos de pepita"
$ws.Cells.Item(1171, 9).Value = 100104005
$ws.Cells.Item(1171, 10).Value = "Pera"
$ws.Cells.Item(1171, 11).Value = "Packham's Triumph"
$ws.Cells.Item(1171, 12).Value = "Primera"
$ws.Cells.Item(1171, 13).Value = 300
$ws.Cells.Item(1171, 14).Value = 14000
$ws.Cells.Item(1171, 15).Value = 14000
$ws.Cells.Item(1171, 16).Value = 14000
$ws.Cells.Item(1171, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(1171, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(1171, 19).Value = 778
$ws.Cells.Item(1171, 20).Value = 18

# Make sure the date cell keeps the same number format as the rest of
# column D (YYYY-MM-DD HH:MM:SS custom format).
$ws.Cells.Item(1171, 4).NumberFormat = $ws.Cells.Item(1172, 4).NumberFormat
